# COREINTL_holdings.xlsx - "Add files via upload" re-upload edit:
#   - refreshed model-holdings-as-of date (2021-05-07 -> 2021-05-10)
#   - refreshed Weight/Percent Change figures for EFA / EEM / Total rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; lift it so the cells below can be written, then
# restore protection afterwards.
$ws.Unprotect()

# Confidentiality footnote text - only the as-of date changed.
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-10 for illustrative purposes only and are subject to change."

# Updated Weight / Percent Change figures.
$ws.Range("D2").Value = 0.8464745141798429
$ws.Range("E2").Value = -0.0053475935828875

$ws.Range("D3").Value = 0.1535254858201571
$ws.Range("E3").Value = -0.01791918083744737

$ws.Range("E4").Value = -0.007277652623674413

$ws.Protect()
